$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 3 ("0002872694" invoice),
# pushing it down to row 6, then fill rows 3-5 with the new invoice data.
$ws.Rows("3:5").Insert()

# Ensure every new cell is formatted as Text before values are written so
# Excel doesn't "helpfully" reinterpret things like "07/08/2025" as a date
# serial or "0003340074" as a number (dropping leading zeros). The shifted
# row (now row 6) keeps its original (default) formatting untouched.
$ws.Range("A3:G5").NumberFormat = "@"

$data = @(
    @("0003340074", "26/08/2025", "19/08/2025", "908,91", "030.663.374", "34191.09123 70053.972934 85972.140009 1 11850000090891", "08.2025 - Energisa - Inst 103340074-8.pdf"),
    @("0000893797", "16/09/2025", "07/08/2025", "0,00", "014.775.807", "74593.10046 27628.019005 01374.441762 7 12060000000000", "08.2025 - Energisa - Inst 8893797-1.pdf"),
    @("893797", "25/09/2025", "18/09/2025", "29,53", "005.598.958", "83650000000-2 29530012000-9 72680172025-1 09200002019-5", "09.2025 - Energisa - Inst 8893797-1 Desligamento.pdf")
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}
